$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (index row=2) updated values
$ws.Range("B2").Value = 0.0067016980621690259
$ws.Range("C2").Value = 0.00077876288205658284
$ws.Range("D2").Value = 0.00021324968540836764
$ws.Range("E2").Value = 0.00011784362499955847

# Row 4 updated values
$ws.Range("B4").Value = 0.01047592140371556
$ws.Range("C4").Value = 0.0072700179959150901
$ws.Range("D4").Value = 0.0091916369931799302
$ws.Range("E4").Value = 0.00099626526775953295

# Row 5 updated values
$ws.Range("B5").Value = 0.0096352374570252142
$ws.Range("C5").Value = 0.018987640380649928
$ws.Range("D5").Value = 0.015372356185406222
$ws.Range("E5").Value = 0.0034478576400140959
